$d = $word.ActiveDocument

# 1) Professional summary paragraph: "all Black and Asian-American voters" -> "50M voters"
$found1 = $d.Content.Find.Execute(
    "Product-focused data scientist with 15+ years building systems that matter. Discovered systematic demographic coding errors affecting all Black and Asian-American voters, developed geospatial ML algorithms improving classification accuracy from 23% to 64%. Expert in translating technical solutions into business value.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Product-focused data scientist with 15+ years building systems that matter. Discovered systematic demographic coding errors affecting 50M voters, developed geospatial ML algorithms improving classification accuracy from 23% to 64%. Expert in translating technical solutions into business value.",
    2)
Write-Output "Summary replaced: $found1"

# 2) Work-experience bullet: same phrase change, but "50M" must become its own bold, colored run
#    (matching the styling already used for the "23%"/"64%" figures in the same sentence)
$bulletRng = $d.Content
$found2 = $bulletRng.Find.Execute(
    "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Discovered systematic race coding errors affecting 50M voters, developed geospatial machine learning algorithms improving demographic classification accuracy from ",
    2)
Write-Output "Bullet text replaced: $found2"

# $bulletRng now spans exactly the text that was just inserted; re-seat its bounds
# explicitly so a fresh Find can run confined to that span, then narrow it to "50M".
$spanStart = $bulletRng.Start
$spanEnd = $bulletRng.End
$bulletRng.Start = $spanStart
$bulletRng.End = $spanEnd

$found3 = $bulletRng.Find.Execute("50M")
Write-Output "50M located within bullet: $found3"
$bulletRng.Font.Bold = 1
$bulletRng.Font.Color = 5258796

# 3) Impact statement paragraph
$found4 = $d.Content.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved electoral prediction accuracy by 22%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved electoral prediction accuracy by 22%",
    2)
Write-Output "Impact replaced: $found4"
